# occurence.txt (first sheet): fix a header typo and collapse the two
# "eventDateBeginning"/"eventDateEnding" Darwin Core columns into a single
# ISO 8601 "eventDate" interval column, per the Darwin Core spec.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Write the merged event-date interval into column I (currently
# "eventDateBeginning") before the now-redundant column J
# ("eventDateEnding") is deleted.
$ws1.Range("I2").Value = "1889-02-01/1889-02-28"
$ws1.Range("I1").Value = "eventDate"

# Fix the "occurenceID" -> "occurrenceID" header typo.
$ws1.Range("C1").Value = "occurrenceID"

# Drop column J ("eventDateEnding"); everything to its right shifts left.
$ws1.Columns.Item(10).Delete()

# The occurence.txt sheet becomes the active tab/selection (it was
# multimedia.txt before), with the cursor parked on the fixed header cell.
$ws1.Activate()
$ws1.Range("C1").Select()
